$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 133, shifting rows 133:144 down to 134:145
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new weekly record
$ws.Range("A133").Value = 10
$ws.Range("B133").Value = "Vega Modelo de Temuco"
$ws.Range("C133").Value = "La Araucanía"
$ws.Range("D133").Value = 45166
$ws.Range("D133").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E133").Value = 9
$ws.Range("F133").Value = 100112010
$ws.Range("G133").Value = "Achicoria"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 100
$ws.Range("K133").Value = 10000
$ws.Range("L133").Value = 10000
$ws.Range("M133").Value = 10000
$ws.Range("N133").Value = "$/caja 18 unidades"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 556
$ws.Range("Q133").Value = 18
$ws.Range("R133").Value = "Hortaliza"
